# Actualización automática 2025-09-08 13:10:08
# Insert a new client row ("BRAVO MANZABA MARIA CECILIA") as row 9 in both the
# "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets (alphabetically between
# "AVILA TORRES RAFAEL ALEJANDRO" and "CARAVEDO PAZMIÑO  JAHAIRA PAMELA"),
# pushing every following client row down by one, and refresh the trailing
# summary row so it reflects the now-larger roster (26 -> 27 clients).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (columns A:R, data rows 2-27, summary row 28)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows.Item(9).Insert()

$ws1.Range("A9").Value = "OFICINA-CATAECSA"
$ws1.Range("B9").Value = "BRAVO MANZABA MARIA CECILIA"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(9, $col).Value = 0
}

# Summary row moved from 28 -> 29; update the "x de 26" -> "x de 27" counters.
$ws1.Range("C29").Value = "0 de 27"
$ws1.Range("D29").Value = "1 de 27"
$ws1.Range("E29").Value = "1 de 27"
$ws1.Range("F29").Value = "0 de 27"
$ws1.Range("G29").Value = "0 de 27"
$ws1.Range("H29").Value = "0 de 27"
$ws1.Range("I29").Value = "0 de 27"
$ws1.Range("J29").Value = "0 de 27"
$ws1.Range("K29").Value = "0 de 27"
$ws1.Range("L29").Value = "2 de 27"
$ws1.Range("M29").Value = "3 de 27"
$ws1.Range("N29").Value = "0 de 27"
$ws1.Range("O29").Value = "0 de 27"
$ws1.Range("P29").Value = "0 de 27"
$ws1.Range("Q29").Value = "0 de 27"
$ws1.Range("R29").Value = "0 de 27"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL" (columns A:G, data rows 2-27, summary row 28)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(9).Insert()

$ws2.Range("A9").Value = "OFICINA-CATAECSA"
$ws2.Range("B9").Value = "BRAVO MANZABA MARIA CECILIA"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(9, $col).Value = 0
}

# Summary row moved from 28 -> 29; totals are unchanged (new row contributes 0).
$ws2.Range("C29").Value = 1076.87
$ws2.Range("D29").Value = 2057.76
$ws2.Range("E29").Value = 1423.94
$ws2.Range("F29").Value = 8670.46
$ws2.Range("G29").Value = 0
